$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Add login screen" -> "Put queries in DML file"          (bullet 2)
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(2)
$p1.Range.Find.Execute("Add login screen", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "Put queries in DML file", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "Put queries in DML file" -> "Put rest of the tables on the
#    website instead of just Database"                        (bullet 3)
#    (scoped to paragraph 3 so bullet 2's brand-new text is untouched)
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(3)
$p2.Range.Find.Execute("Put queries in DML file", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "Put rest of the tables on the website instead of just Database", 2) | Out-Null

# ------------------------------------------------------------------
# 3) "Take out league name from team table in sample data" moves from
#    outline level 0 to level 1 (indented under the previous bullet)
# ------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$p7.Range.ListFormat.ListIndent()

# ------------------------------------------------------------------
# 4) Insert a blank paragraph, then a "Questions" paragraph (plain,
#    no list numbering) right after the "Take out league name..." bullet
# ------------------------------------------------------------------
$p7.Range.InsertParagraphAfter()
$blank = $d.Paragraphs.Item(8)
$blank.Range.ListFormat.RemoveNumbers()
$blank.Range.Style = $d.Styles.Item("Normal")

$blank.Range.InsertParagraphAfter()
$questions = $d.Paragraphs.Item(9)
$questions.Range.Text = "Questions"

# ------------------------------------------------------------------
# 5) The old "Fix auto increment in php file" paragraph becomes the
#    start of a brand-new numbered list: "Login screen"
# ------------------------------------------------------------------
$oldLast = $d.Paragraphs.Item(10)
$bodyRange = $d.Range($oldLast.Range.Start, $oldLast.Range.End - 1)
$bodyRange.Text = "Login screen"

$loginScreen = $d.Paragraphs.Item(10).Range
$loginScreen.ListFormat.RemoveNumbers()
$loginScreen.ListFormat.ApplyNumberDefault()

# ------------------------------------------------------------------
# 6) Append the remaining "Questions" bullets, reusing numId 2 by
#    chaining InsertParagraphAfter from the "Login screen" paragraph
# ------------------------------------------------------------------
$d.Paragraphs.Item(10).Range.InsertParagraphAfter()
$usersAdmin = $d.Paragraphs.Item(11)
$usersAdmin.Range.Text = "Do we need users and admin"
$usersAdmin.Range.ListFormat.ListIndent()

$d.Paragraphs.Item(11).Range.InsertParagraphAfter()
$signUp = $d.Paragraphs.Item(12)
$signUp.Range.Text = "Do we need a " + [char]0x201C + "sign up" + [char]0x201D + " thing for new users or is one username and password fine"

$d.Paragraphs.Item(12).Range.InsertParagraphAfter()
$columns = $d.Paragraphs.Item(13)
$columns.Range.Text = "How to see if columns have same value"
$columns.Range.ListFormat.ListOutdent()

$d.Paragraphs.Item(13).Range.InsertParagraphAfter()
$dmlFile = $d.Paragraphs.Item(14)
$dmlFile.Range.Text = "Do we need a DML file "
$dmlFileP = $d.Paragraphs.Item(14)
$tail = $d.Range($dmlFileP.Range.End - 1, $dmlFileP.Range.End - 1)
$tail.InsertAfter("of our 20 queries")
# force the appended text to stay in its own run instead of being
# merged back into the previous one
$secondRun = $d.Range($d.Paragraphs.Item(14).Range.End - 1 - ("of our 20 queries").Length, $d.Paragraphs.Item(14).Range.End - 1)
$secondRun.Bold = 1
$secondRun.Bold = 0

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
